# Apply the Sem1_Plan.xlsx update:
#  - Insert a new Gantt task "Analyse questionnaire results" between
#    "Create and release Lua questionnaire" and "Create a list of requirements".
#  - Push out the "Research" phase end dates to accommodate the new task.
#  - Pull the "Select an appropriate compiler" task earlier.
#  - Keep the Min/Max summary formulas and the chart pointed at the right ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt Chart")

# --- Insert the new row for the task, shifting everything below down ---
$ws.Rows(10).Insert()

# Copy the formatting of a normal task row (now row 12, "Create a list of
# requirements") onto the freshly inserted (blank) row 10.
$ws.Range("B12:E12").Copy()
$ws.Range("B10:E10").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Research phase: push End Date out for rows 7-9 ---
$ws.Cells.Item(7, 4).Value = 44162
$ws.Cells.Item(8, 4).Value = 44168
$ws.Cells.Item(9, 4).Value = 44168

# --- New task row 10: "Analyse questionnaire results" ---
$ws.Cells.Item(10, 2).Value = "Analyse questionnaire results"
$ws.Cells.Item(10, 3).Value = 44168
$ws.Cells.Item(10, 4).Value = 44173
$ws.Cells.Item(10, 5).Formula = "=D10-C10"

# --- "Select an appropriate compiler" (now row 13) moves earlier ---
$ws.Cells.Item(13, 3).Value = 44162
$ws.Cells.Item(13, 4).Value = 44173

# Re-enter E12 (now standalone) and E13 (now the head of a shared
# E13:E15 duration formula block) to mirror the author's edit.
$ws.Cells.Item(12, 5).Formula = "=D12-C12"
$ws.Cells.Item(13, 5).Formula = "=D13-C13"
$ws.Cells.Item(14, 5).Formula = "=D14-C14"
$ws.Cells.Item(15, 5).Formula = "=D15-C15"

$ws.Range("D18").Select() | Out-Null
